# FALECPV-CajaChica.xlsx
# "reporte de caja chica con el estado" -> add an ESTADO (status) column to
# the petty-cash report table, right after "DOC REFRENCIA" (column E) and
# before "NOTAS" (previously column F, now shifted to G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; this shifts the existing NOTAS / REPOSICIÓN /
# GASTO / USUARIO / REGISTRO columns one place to the right (F->G, G->H,
# H->I, I->J, J->K) and extends the title merge (A2:J2 -> A2:K2) and the
# sheet dimension (A2:J10 -> A2:K10) automatically.
$null = $ws.Columns("F:F").Insert()

# The new column inherits the width of its neighbour to the left (column
# E, "DOC REFRENCIA" = 22.5 characters wide).
$ws.Columns("F:F").ColumnWidth = 21.67

# Header for the new column in the table row.
$ws.Range("F10").Value = "ESTADO"

# Leave the selection on the newly added header cell, same as the author.
$null = $ws.Range("F10").Select()
